$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Trim type" column (E) holds "none" for every data row (E2:E10).
# Change it to "both" for the whole output-file metadata table.
$ws.Range("E2:E10").Value = "both"

# Reflect the edit in the sheet's current selection (active cell E2,
# selected range E2:E10), matching the post-edit view state.
$ws.Range("E2:E10").Select()
